$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Laporan penjualan"

$ws.Range("A1").Value = "Kuartal Tahun"
$ws.Range("A1").Characters(1,12).Font.Color = 0
$ws.Range("A1").Characters(1,12).Font.Name = "Aptos Narrow"
$ws.Range("A1").Characters(1,12).Font.Size = 11
$ws.Range("A1").Characters(13,1).Font.Color = 0
$ws.Range("A1").Characters(13,1).Font.Name = "Aptos Narrow"
$ws.Range("A1").Characters(13,1).Font.Size = 11

$ws.Range("B1").Value = "Barat Tengah"
$ws.Range("B1").Characters(1,11).Font.Color = 0
$ws.Range("B1").Characters(1,11).Font.Name = "Aptos Narrow"
$ws.Range("B1").Characters(1,11).Font.Size = 11
$ws.Range("B1").Characters(12,1).Font.Color = 0
$ws.Range("B1").Characters(12,1).Font.Name = "Aptos Narrow"
$ws.Range("B1").Characters(12,1).Font.Size = 11

$ws.Range("C1").Value = "Gunung"
$ws.Range("C1").Characters(1,5).Font.Color = 0
$ws.Range("C1").Characters(1,5).Font.Name = "Aptos Narrow"
$ws.Range("C1").Characters(1,5).Font.Size = 11
$ws.Range("C1").Characters(6,1).Font.Color = 0
$ws.Range("C1").Characters(6,1).Font.Name = "Aptos Narrow"
$ws.Range("C1").Characters(6,1).Font.Size = 11

$ws.Range("D1").Value = "Timur laut"
$ws.Range("D1").Characters(1,9).Font.Color = 0
$ws.Range("D1").Characters(1,9).Font.Name = "Aptos Narrow"
$ws.Range("D1").Characters(1,9).Font.Size = 11
$ws.Range("D1").Characters(10,1).Font.Color = 0
$ws.Range("D1").Characters(10,1).Font.Name = "Aptos Narrow"
$ws.Range("D1").Characters(10,1).Font.Size = 11

$ws.Range("E1").Value = "Selatan"
$ws.Range("E1").Characters(1,6).Font.Color = 0
$ws.Range("E1").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("E1").Characters(1,6).Font.Size = 11
$ws.Range("E1").Characters(7,1).Font.Color = 0
$ws.Range("E1").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("E1").Characters(7,1).Font.Size = 11

$ws.Range("F1").Value = "Tenggara"
$ws.Range("F1").Characters(1,7).Font.Color = 0
$ws.Range("F1").Characters(1,7).Font.Name = "Aptos Narrow"
$ws.Range("F1").Characters(1,7).Font.Size = 11
$ws.Range("F1").Characters(8,1).Font.Color = 0
$ws.Range("F1").Characters(8,1).Font.Name = "Aptos Narrow"
$ws.Range("F1").Characters(8,1).Font.Size = 11

$ws.Range("G1").Value = "Barat"
$ws.Range("G1").Characters(1,4).Font.Color = 0
$ws.Range("G1").Characters(1,4).Font.Name = "Aptos Narrow"
$ws.Range("G1").Characters(1,4).Font.Size = 11
$ws.Range("G1").Characters(5,1).Font.Color = 0
$ws.Range("G1").Characters(5,1).Font.Name = "Aptos Narrow"
$ws.Range("G1").Characters(5,1).Font.Size = 11

$ws.Range("A2").Value = "Q1 2022"
$ws.Range("A2").Characters(1,6).Font.Color = 0
$ws.Range("A2").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A2").Characters(1,6).Font.Size = 11
$ws.Range("A2").Characters(7,1).Font.Color = 0
$ws.Range("A2").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A2").Characters(7,1).Font.Size = 11

$ws.Range("A3").Value = "Q2 2022"
$ws.Range("A3").Characters(1,6).Font.Color = 0
$ws.Range("A3").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A3").Characters(1,6).Font.Size = 11
$ws.Range("A3").Characters(7,1).Font.Color = 0
$ws.Range("A3").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A3").Characters(7,1).Font.Size = 11

$ws.Range("A4").Value = "Q3 2022"
$ws.Range("A4").Characters(1,6).Font.Color = 0
$ws.Range("A4").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A4").Characters(1,6).Font.Size = 11
$ws.Range("A4").Characters(7,1).Font.Color = 0
$ws.Range("A4").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A4").Characters(7,1).Font.Size = 11

$ws.Range("A5").Value = "Q4 2022"
$ws.Range("A5").Characters(1,6).Font.Color = 0
$ws.Range("A5").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A5").Characters(1,6).Font.Size = 11
$ws.Range("A5").Characters(7,1).Font.Color = 0
$ws.Range("A5").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A5").Characters(7,1).Font.Size = 11

$ws.Range("A6").Value = "Q1-2023"
$ws.Range("A6").Characters(1,6).Font.Color = 0
$ws.Range("A6").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A6").Characters(1,6).Font.Size = 11
$ws.Range("A6").Characters(7,1).Font.Color = 0
$ws.Range("A6").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A6").Characters(7,1).Font.Size = 11

$ws.Range("A7").Value = "Q2-2023"
$ws.Range("A7").Characters(1,6).Font.Color = 0
$ws.Range("A7").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A7").Characters(1,6).Font.Size = 11
$ws.Range("A7").Characters(7,1).Font.Color = 0
$ws.Range("A7").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A7").Characters(7,1).Font.Size = 11

$ws.Range("A8").Value = "Q3-2023"
$ws.Range("A8").Characters(1,6).Font.Color = 0
$ws.Range("A8").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A8").Characters(1,6).Font.Size = 11
$ws.Range("A8").Characters(7,1).Font.Color = 0
$ws.Range("A8").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A8").Characters(7,1).Font.Size = 11

$ws.Range("A9").Value = "Q4-2023"
$ws.Range("A9").Characters(1,6).Font.Color = 0
$ws.Range("A9").Characters(1,6).Font.Name = "Aptos Narrow"
$ws.Range("A9").Characters(1,6).Font.Size = 11
$ws.Range("A9").Characters(7,1).Font.Color = 0
$ws.Range("A9").Characters(7,1).Font.Name = "Aptos Narrow"
$ws.Range("A9").Characters(7,1).Font.Size = 11
